# CIERRE 11 AGO 23
# Update the "VALES DE INSENTIVOS" vale: amount 5000 -> 6000 (CINCO -> SEIS
# MIL PESOS) and month JUNIO -> JULIO, then leave that sheet selected
# (cursor on J9) as the active tab instead of "ARQUITECTO        ".

$wb = $excel.ActiveWorkbook

$wsArquitecto = $wb.Worksheets.Item("ARQUITECTO        ")
$wsVales      = $wb.Worksheets.Item("VALES DE INSENTIVOS")

# --- Update the incentive vale amount and wording ---
$wsVales.Range("D1").Value = 6000
$wsVales.Range("A2").Value = "SEIS   MIL     PESOS 00/100 M.N."
$wsVales.Range("A4").Value = "INCENTIVO DEL MES DE  JULIO     2023"

# --- Switch the active/selected sheet + cursor position ---
[void]$wsVales.Activate()
[void]$wsVales.Range("J9").Select()
